$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete rows (old rows 5, 6 and 7) so the table only keeps
# the first three data rows (new rows 2-4).
$ws.Rows("5:7").Delete()

# --- Row 2 (ECs -> FAPs, unchanged labels) -------------------------------
$ws.Range("G2").Value = 2.020378666666666
$ws.Range("H2").Value = 6.061135999999999
$ws.Range("I2").Value = 0.2600911804892308
$ws.Range("J2").Value = 0.2600911804892308
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.4293345074239999
$ws.Range("R2").Value = 3.864010566816
$ws.Range("S2").Value = 0.2600911804892308
$ws.Range("T2").Value = 0.2600911804892308

# --- Row 3 (ECs -> FAPs, MuSCs -> FAPs) ----------------------------------
$ws.Range("A3").Value = "FAPs"
$ws.Range("D3").Value = "FAPs"
$ws.Range("G3").Value = 3.623135666666666
$ws.Range("H3").Value = 10.869407
$ws.Range("I3").Value = 0.4664203043534923
$ws.Range("J3").Value = 0.4664203043534922
$ws.Range("M3").Value = 0.212502
$ws.Range("N3").Value = 0.637506
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.7699235754379999
$ws.Range("R3").Value = 6.929312178941999
$ws.Range("S3").Value = 0.4664203043534923
$ws.Range("T3").Value = 0.4664203043534922

# --- Row 4 (FAPs -> MuSCs, target cluster unchanged) ---------------------
$ws.Range("A4").Value = "MuSCs"
$ws.Range("G4").Value = 2.124448666666666
$ws.Range("H4").Value = 6.373346
$ws.Range("I4").Value = 0.273488515157277
$ws.Range("J4").Value = 0.273488515157277
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.451449590564
$ws.Range("R4").Value = 4.063046315076
$ws.Range("S4").Value = 0.273488515157277
$ws.Range("T4").Value = 0.273488515157277
